# Fruta / hortaliza, semanal
# Insert two new weekly price rows (16:17) for "Cilantro" at
# "Terminal Hortofrutícola Agro Chillán", shifting the existing data
# rows down by two (old row 16 -> 18, ..., old row 58 -> 60).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 16-58 down to 18-60, inserting two blank rows.
$ws.Rows("16:17").Insert()

# New row 16: Primera, 2022-07-22 (serial 44764)
$ws.Range("A16").Value = 7
$ws.Range("B16").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C16").Value = "Ñuble"
$ws.Range("D16").Value = 44764
$ws.Range("E16").Value = 16
$ws.Range("F16").Value = 100112040
$ws.Range("G16").Value = "Cilantro"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 160
$ws.Range("K16").Value = 700
$ws.Range("L16").Value = 800
$ws.Range("M16").Value = 750
$ws.Range("N16").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O16").Value = "Provincia de Diguillín"
$ws.Range("P16").Value = 750
$ws.Range("Q16").Value = 1
$ws.Range("R16").Value = "Hortaliza"

# New row 17: Segunda, 2022-07-22 (serial 44764)
$ws.Range("A17").Value = 7
$ws.Range("B17").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C17").Value = "Ñuble"
$ws.Range("D17").Value = 44764
$ws.Range("E17").Value = 16
$ws.Range("F17").Value = 100112040
$ws.Range("G17").Value = "Cilantro"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Segunda"
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 600
$ws.Range("L17").Value = 600
$ws.Range("M17").Value = 600
$ws.Range("N17").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O17").Value = "Provincia de Diguillín"
$ws.Range("P17").Value = 600
$ws.Range("Q17").Value = 1
$ws.Range("R17").Value = "Hortaliza"
